$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.772.35'
$ws.Range('E2').Value = '  -1.41%  '
$ws.Range('D3').Value = '2.594.64'
$ws.Range('E3').Value = '  -2.06%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = "'551.68"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.59%  '
$ws.Range('D6').Value = "'143.03"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.19%  '
$ws.Range('D8').Value = "'0.603"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +5.27%  '
$ws.Range('D9').Value = "'6.77"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.59%  '
$ws.Range('E10').Value = '  -1.92%  '
$ws.Range('E11').Value = '  +4.97%  '
$ws.Range('D12').Value = "'0.336"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.82%  '
$ws.Range('D13').Value = '3.055.18'
$ws.Range('E13').Value = '  -1.92%  '
$ws.Range('D14').Value = '58.730.52'
$ws.Range('E14').Value = '  -1.34%  '
$ws.Range('D15').Value = "'20.85"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.13%  '
$ws.Range('D16').Value = '2.597.55'
$ws.Range('E16').Value = '  -1.75%  '
$ws.Range('E17').Value = '  -2.20%  '
$ws.Range('E18').Value = '  +1.23%  '
$ws.Range('D19').Value = "'337.17"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.88%  '
$ws.Range('D20').Value = "'10.05"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.73%  '
$ws.Range('E21').Value = '  -1.07%  '
$ws.Range('D23').Value = "'66.82"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.27%  '
$ws.Range('D24').Value = "'0.427"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.36%  '
$ws.Range('D25').Value = "'0.999"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('D26').Value = "'0.159"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.31%  '
$ws.Range('D27').Value = "'7.13"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.16%  '
$ws.Range('D28').Value = '0.0₃0752'
$ws.Range('E28').Value = '  +0.77%  '
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('E30').Value = '  +1.49%  '
$ws.Range('D32').Value = "'154.66"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.49%  '
$ws.Range('D33').Value = "'18.94"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.12%  '
$ws.Range('D34').Value = "'3.92"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.04%  '
$ws.Range('D35').Value = "'0.893"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.14%  '
$ws.Range('D36').Value = "'1.12"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.23%  '
$ws.Range('E38').Value = '  +1.16%  '
$ws.Range('D39').Value = "'0.831"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.78%  '
$ws.Range('D40').Value = "'3.60"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.14%  '
$ws.Range('D41').Value = "'283.60"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.02%  '
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('D43').Value = "'0.598"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.45%  '
$ws.Range('D44').Value = "'0.0958"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.31%  '
$ws.Range('E45').Value = '  -1.20%  '
$ws.Range('D46').Value = "'0.0533"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.21%  '
$ws.Range('E47').Value = '  -0.12%  '
$ws.Range('D48').Value = '1.946.63'
$ws.Range('E48').Value = '  -1.06%  '
$ws.Range('D49').Value = "'118.42"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +6.58%  '
$ws.Range('D50').Value = "'17.85"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.62%  '
$ws.Range('D51').Value = "'4.41"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.48%  '
